# Apply the authored changes described in the commit:
#  - Refresh the auto date field placeholders (11/30/2022 -> 12/1/2022)
#    on the slide master and every slide layout.
#  - Shorten the subtitle name on slide 1.
#  - Fix the "Front-end" heading spacing on slide 10.
#  - Shorten the title on slide 4.
#  - Nudge a picture's horizontal position on slide 9.

$p = $ppt.ActivePresentation

$oldDate = "11/30/2022"
$newDate = "12/1/2022"

# --- Slide Master date placeholder ---
$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        if ($sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Every slide layout's date placeholder ---
for ($j = 1; $j -le $m.CustomLayouts.Count; $j++) {
    $lay = $m.CustomLayouts.Item($j)
    for ($i = 1; $i -le $lay.Shapes.Count; $i++) {
        $sh = $lay.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Slide 1: shorten the presenter's name ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "Silvio Tavares"

# --- Slide 4: shorten the title ---
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Amizades"

# --- Slide 9: nudge the screenshot image to the right ---
# Shape.Left is expressed in points (1 pt = 12700 EMU); the target
# offset is 1067107 EMU. A tiny epsilon keeps the float-to-EMU
# rounding inside this runtime from landing one EMU short.
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(2).Left = (1067107 / 12700) + 0.00002

# --- Slide 10: "Front-" -> "Front " (keep the following "end" run intact) ---
$s10 = $p.Slides.Item(10)
$tr10 = $s10.Shapes.Item(2).TextFrame.TextRange
$para1 = $tr10.Paragraphs(1, 1)
$firstRun = $para1.Characters(1, 6)
$firstRun.Text = "Front "
